# VAN-1811: Prepare and write FUNCTIONAL test cases and test scripts
#
# Updates the "Order_Assign_ByPreviousDoc" test-data row (row 2 of Sheet1)
# with new test values:
#   - AV2 ("FuncLoc")     : ASMPD503ALVINCOMMUN -> ABCD690970
#   - AX2 ("Previous Doc"): 2152430001           -> 9605206304

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AV2").Value = "ABCD690970"
$ws.Range("AX2").Value = "9605206304"
